# Auto-generated edit script: updates crypto price/volume table cells
# to match the refreshed GitHub Actions data pull, and fixes the
# THORChain / LidoDAOToken row order (rows 36-37).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $text) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "38.638.64"
Set-TextCell "E2" "  +2.36%  "
Set-TextCell "D3" "2.104.85"
Set-TextCell "E3" "  +3.80%  "
Set-TextCell "E4" "  -0.02%  "
Set-TextCell "D5" "229.64"
Set-TextCell "E5" "  +1.03%  "
Set-TextCell "D6" "0.616"
Set-TextCell "E6" "  +1.35%  "
Set-TextCell "D7" "61.55"
Set-TextCell "E7" "  +2.89%  "
Set-TextCell "D8" "1.00"
Set-TextCell "E8" "  -0.01%  "
Set-TextCell "D9" "0.382"
Set-TextCell "E9" "  +1.72%  "
Set-TextCell "D10" "0.0847"
Set-TextCell "E10" "  +3.30%  "
Set-TextCell "E11" "  +0.51%  "
Set-TextCell "D12" "2.415.52"
Set-TextCell "E12" "  +3.80%  "
Set-TextCell "D13" "14.84"
Set-TextCell "E13" "  +2.86%  "
Set-TextCell "D14" "22.41"
Set-TextCell "E14" "  +6.56%  "
Set-TextCell "D15" "0.783"
Set-TextCell "E15" "  +1.72%  "
Set-TextCell "E16" "  +5.73%  "
Set-TextCell "D17" "2.103.05"
Set-TextCell "E17" "  +4.10%  "
Set-TextCell "D18" "38.525.54"
Set-TextCell "E18" "  +2.26%  "
Set-TextCell "D19" "6.03"
Set-TextCell "E19" "  +2.33%  "
Set-TextCell "D20" "70.58"
Set-TextCell "E20" "  +1.68%  "
Set-TextCell "D21" "0.0₃0837"
Set-TextCell "E21" "  +1.79%  "
Set-TextCell "D22" "226.70"
Set-TextCell "E22" "  +1.31%  "
Set-TextCell "E23" "  -0.10%  "
Set-TextCell "D24" "2.44"
Set-TextCell "E24" "  +3.08%  "
Set-TextCell "D25" "2.32"
Set-TextCell "E25" "  +3.60%  "
Set-TextCell "D26" "169.82"
Set-TextCell "E26" "  +1.19%  "
Set-TextCell "D27" "9.44"
Set-TextCell "E27" "  +1.27%  "
Set-TextCell "D28" "0.131"
Set-TextCell "E28" "  +2.18%  "
Set-TextCell "D29" "19.10"
Set-TextCell "E29" "  +1.75%  "
Set-TextCell "D30" "1.37"
Set-TextCell "E30" "  +9.08%  "
Set-TextCell "E31" "  +0.11%  "
Set-TextCell "D32" "2.34"
Set-TextCell "E32" "  +5.87%  "
Set-TextCell "E33" "  +6.19%  "
Set-TextCell "D34" "4.48"
Set-TextCell "E34" "  +2.44%  "
Set-TextCell "D35" "0.0606"
Set-TextCell "E35" "  +0.36%  "
Set-TextCell "B36" "THORChain"
Set-TextCell "C36" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell "D36" "6.48"
Set-TextCell "E36" "  +1.20%  "
Set-TextCell "B37" "LidoDAOToken"
Set-TextCell "C37" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D37" "2.40"
Set-TextCell "E37" "  +4.07%  "
Set-TextCell "D38" "3.51"
Set-TextCell "E38" "  +3.16%  "
Set-TextCell "D39" "0.999"
Set-TextCell "E39" "  -0.15%  "
Set-TextCell "D40" "18.40"
Set-TextCell "E40" "  +2.91%  "
Set-TextCell "D41" "1.539.40"
Set-TextCell "E41" "  +0.36%  "
Set-TextCell "D42" "100.20"
Set-TextCell "E42" "  +4.82%  "
Set-TextCell "D43" "0.0221"
Set-TextCell "E43" "  +2.49%  "
Set-TextCell "E44" "  +0.86%  "
Set-TextCell "D45" "0.0913"
Set-TextCell "E45" "  +0.62%  "
Set-TextCell "D46" "4.18"
Set-TextCell "E46" "  +2.97%  "
Set-TextCell "D47" "1.12"
Set-TextCell "E47" "  +1.50%  "
Set-TextCell "D48" "7.55"
Set-TextCell "E48" "  +6.29%  "
Set-TextCell "E49" "  +3.88%  "
Set-TextCell "D50" "2.99"
Set-TextCell "E50" "  +0.90%  "
Set-TextCell "D51" "2.301.41"
Set-TextCell "E51" "  +3.81%  "
